$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 397.5
$ws.Range("B3").Value = 359.4
$ws.Range("C3").Value = 405.1
$ws.Range("C4").Value = 405.6
$ws.Range("C11").Value = 359
$ws.Range("C18").Value = 429.6
